# Update the Levine/Toilolo TE 2018 sheet with new "height" and "weight" columns.
# The existing "fantasy points" column (E) is relabeled to "height", a new
# "weight" column is inserted in F, and the original "fantasy points" values
# move out into the new column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 17

# Capture the existing "fantasy points" values (column E, rows 2-17) before
# they get overwritten, so they can be moved into the new column G.
$fantasyPoints = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $fantasyPoints[$r] = $ws.Cells.Item($r, 5).Value()
}

# Re-label column E's header from "fantasy points" to "height", and add the
# new "weight" / "fantasy points" headers in F1 / G1, matching the header
# formatting already used by columns B-D.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:G1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the new data: column E = height (constant, in feet), column F =
# weight (constant, in lbs), column G = the original fantasy points values.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.666666666666667
    $ws.Cells.Item($r, 6).Value = 268
    $ws.Cells.Item($r, 7).Value = $fantasyPoints[$r]
}
